$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textFormatCells = @("D5", "D7", "D10", "D13", "D14", "D16", "D20", "D22", "D26", "D28", "D29", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D42", "D43", "D44", "D46", "D49", "D50")
foreach ($c in $textFormatCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = '38.083.59'
$ws.Range("E2").Value = '  +0.59%  '

$ws.Range("D3").Value = '2.092.49'
$ws.Range("E3").Value = '  +2.97%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '228.77'
$ws.Range("E5").Value = '  +0.44%  '

$ws.Range("E6").Value = '  +0.20%  '

$ws.Range("D7").Value = '60.92'
$ws.Range("E7").Value = '  +0.94%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("E9").Value = '  -0.16%  '

$ws.Range("D10").Value = '0.0841'
$ws.Range("E10").Value = '  +2.87%  '

$ws.Range("E11").Value = '  +0.02%  '

$ws.Range("D12").Value = '2.400.65'
$ws.Range("E12").Value = '  +2.87%  '

$ws.Range("D13").Value = '14.66'
$ws.Range("E13").Value = '  +0.81%  '

$ws.Range("D14").Value = '22.12'
$ws.Range("E14").Value = '  +3.32%  '

$ws.Range("E15").Value = '  +6.53%  '

$ws.Range("D16").Value = '0.775'
$ws.Range("E16").Value = '  +1.88%  '

$ws.Range("D17").Value = '2.088.57'
$ws.Range("E17").Value = '  +1.92%  '

$ws.Range("D18").Value = '37.695.66'
$ws.Range("E18").Value = '  -0.30%  '

$ws.Range("E19").Value = '  +2.07%  '

$ws.Range("D20").Value = '70.04'
$ws.Range("E20").Value = '  +0.28%  '

$ws.Range("D21").Value = '0.0₃0840'
$ws.Range("E21").Value = '  +1.57%  '

$ws.Range("D22").Value = '223.70'
$ws.Range("E22").Value = '  -0.33%  '

$ws.Range("E23").Value = '  +0.53%  '

$ws.Range("E24").Value = '  +0.20%  '

$ws.Range("E25").Value = '  +3.16%  '

$ws.Range("D26").Value = '169.68'
$ws.Range("E26").Value = '  +1.52%  '

$ws.Range("E27").Value = '  +0.93%  '

$ws.Range("D28").Value = '0.132'
$ws.Range("E28").Value = '  +3.48%  '

$ws.Range("D29").Value = '18.98'
$ws.Range("E29").Value = '  +0.43%  '

$ws.Range("E30").Value = '  +3.97%  '

$ws.Range("E31").Value = '  -0.28%  '

$ws.Range("E32").Value = '  +10.37%  '

$ws.Range("D33").Value = '4.44'
$ws.Range("E33").Value = '  +0.85%  '

$ws.Range("D34").Value = '4.67'
$ws.Range("E34").Value = '  +3.61%  '

$ws.Range("D35").Value = '0.0607'
$ws.Range("E35").Value = '  +0.21%  '

$ws.Range("D36").Value = '2.42'
$ws.Range("E36").Value = '  +6.24%  '

$ws.Range("E37").Value = '  -0.15%  '

$ws.Range("D38").Value = '3.54'
$ws.Range("E38").Value = '  +8.72%  '

$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  +0.01%  '

$ws.Range("D40").Value = '18.00'
$ws.Range("E40").Value = '  +5.09%  '

$ws.Range("D41").Value = '1.547.36'
$ws.Range("E41").Value = '  +1.56%  '

$ws.Range("D42").Value = '99.97'
$ws.Range("E42").Value = '  +4.05%  '

$ws.Range("D43").Value = '0.0219'
$ws.Range("E43").Value = '  +0.78%  '

$ws.Range("D44").Value = '2.82'
$ws.Range("E44").Value = '  -0.21%  '

$ws.Range("E45").Value = '  -0.93%  '

$ws.Range("D46").Value = '4.13'
$ws.Range("E46").Value = '  +3.34%  '

$ws.Range("E47").Value = '  +0.67%  '

$ws.Range("E48").Value = '  +1.51%  '

$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").Value = '7.23'
$ws.Range("E49").Value = '  +1.77%  '

$ws.Range("B50").Value = 'MXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D50").Value = '3.00'
$ws.Range("E50").Value = '  +1.25%  '

$ws.Range("D51").Value = '2.287.12'
$ws.Range("E51").Value = '  +2.91%  '
